# ============================================================================
# feat: add 2022-Q3 data
#
# 1) Insert a new "2022-Q3" worksheet right after "总计", populated with the
#    fund-holding breakdown for the quarter (copied from the "2022-Q1" sheet
#    so formatting/styles line up, then all values overwritten).
# 2) Insert a new leading row in "总计" for the "2022-Q3" totals and
#    renumber the trailing index column.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Step 1: "总计" (totals) sheet — insert the new 2022-Q3 row at the top.
# ----------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# The freshly inserted row lost column A's number style (centered/bordered) —
# borrow it back from the row below, then strip the stray partial style
# Excel's row-insert gave columns B:D.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 29
$summary.Range("D2").Value = 26.95

# Rows 3..8 kept their old index-column values after the insert shifted them
# down one row; renumber so the 0-based counter matches row order again.
for ($r = 3; $r -le 8; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ----------------------------------------------------------------------------
# Step 2: New "2022-Q3" worksheet with the fund holdings table.
# ----------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q1")
$template.Copy($null, $wb.Worksheets.Item("总计"))
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q3"

# The template only has 27 data rows (rows 2-28); this quarter has 29 data
# rows (rows 2-30), so stamp the extra two rows' formatting from row 28.
$ws.Range("A28:H28").Copy()
$ws.Range("A29:H30").PasteSpecial(-4122)

# Force columns D-G (numeric-looking text: fund size/position/value figures)
# to store plain text so things like the trailing zero in "91.10" survive
# instead of being coerced to a number — matches the source data's types.
# Column B (fund code, e.g. "008545") also needs this to keep leading zeros.
# Column C (fund name) is never numeric-looking, so it's left alone.
$ws.Range("B2:B30").NumberFormat = "@"
$ws.Range("D2:G30").NumberFormat = "@"

$rows = @(
    "0`t008545`t泓德丰润三年持有期混合`t74.71`t88.42`t7.37`t5.5061`t3",
    "1`t011058`t景顺长城成长龙头一年持有期混合A`t51.58`t91.71`t6.40`t3.3011`t1",
    "2`t260101`t景顺长城优选混合`t49.69`t76.36`t5.43`t2.6982`t2",
    "3`t001975`t景顺长城环保优势股票`t44.53`t92.46`t5.71`t2.5427`t3",
    "4`t006435`t景顺长城创新成长混合`t36.71`t92.71`t5.90`t2.1659`t2",
    "5`t005395`t泓德臻远回报灵活配置混合`t29.76`t93.35`t6.68`t1.9880`t4",
    "6`t010864`t泓德卓远混合A`t22.84`t92.87`t6.48`t1.4800`t4",
    "7`t001500`t泓德远见回报混合`t21.21`t92.75`t6.75`t1.4317`t8",
    "8`t010003`t景顺长城电子信息产业股票A`t18.10`t92.74`t5.85`t1.0588`t1",
    "9`t004965`t泓德致远混合A`t16.19`t46.90`t5.80`t0.9390`t3",
    "10`t009376`t景顺长城成长领航混合`t13.15`t92.87`t5.83`t0.7666`t2",
    "11`t010865`t泓德卓远混合C`t10.32`t92.87`t6.48`t0.6687`t4",
    "12`t010004`t景顺长城电子信息产业股票C`t7.29`t92.74`t5.85`t0.4265`t1",
    "13`t011059`t景顺长城成长龙头一年持有期混合C`t6.31`t91.71`t6.40`t0.4038`t1",
    "14`t003501`t泰达宏利睿智稳健灵活配置混合A`t10.30`t76.72`t3.29`t0.3389`t1",
    "15`t162204`t泰达宏利行业精选混合A`t8.57`t81.08`t3.43`t0.2940`t2",
    "16`t013280`t泰达宏利睿智稳健灵活配置混合C`t5.95`t76.72`t3.29`t0.1958`t1",
    "17`t260111`t景顺长城公司治理混合`t3.52`t91.70`t5.08`t0.1788`t2",
    "18`t513360`t博时中证全球中国教育主题ETF（QDII）`t4.81`t99.43`t3.38`t0.1626`t8",
    "19`t004966`t泓德致远混合C`t2.54`t46.90`t5.80`t0.1473`t3",
    "20`t015601`t泰达宏利行业精选混合C`t4.13`t81.08`t3.43`t0.1417`t2",
    "21`t012320`t西藏东财国证消费电子主题指数增强C`t0.65`t94.49`t3.71`t0.0241`t7",
    "22`t159725`t工银瑞信中证线上消费主题ETF`t0.57`t98.42`t3.87`t0.0221`t6",
    "23`t002244`t景顺长城低碳科技主题灵活配置混合`t0.56`t62.99`t3.92`t0.0220`t3",
    "24`t012319`t西藏东财国证消费电子主题指数增强A`t0.37`t94.49`t3.71`t0.0137`t7",
    "25`t159728`t南方国证在线消费ETF`t0.33`t99.99`t3.89`t0.0128`t5",
    "26`t001535`t景顺长城改革机遇灵活配置混合A`t0.26`t64.66`t3.09`t0.0080`t5",
    "27`t410010`t华富中小企业100指数增强`t0.18`t93.66`t3.09`t0.0056`t10",
    "28`t007945`t景顺长城改革机遇灵活配置混合C`t0.06`t64.66`t3.09`t0.0019`t5"
)

$r = 2
foreach ($line in $rows) {
    $f = $line.Split("`t")
    $ws.Cells.Item($r, 1).Value = [int]$f[0]
    $ws.Cells.Item($r, 2).Value = $f[1]
    $ws.Cells.Item($r, 3).Value = $f[2]
    $ws.Cells.Item($r, 4).Value = $f[3]
    $ws.Cells.Item($r, 5).Value = $f[4]
    $ws.Cells.Item($r, 6).Value = $f[5]
    $ws.Cells.Item($r, 7).Value = $f[6]
    $ws.Cells.Item($r, 8).Value = [int]$f[7]
    $r = $r + 1
}
